$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 111866252
$ws.Range("B3").Value = 78228
$ws.Range("E3").Value = 6453
$ws.Range("F3").Value = "Vedskivlav"
$ws.Range("G3").Value = "Hertelidea botryosa"
$ws.Range("H3").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q3").Value = 702681
$ws.Range("R3").Value = 7299925
$ws.Range("A4").Value = 111866194
$ws.Range("B4").Value = 90816
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 2059
$ws.Range("F4").Value = "Skrovlig taggsvamp"
$ws.Range("G4").Value = "Hydnellum scabrosum"
$ws.Range("H4").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q4").Value = 702687
$ws.Range("R4").Value = 7299920
$ws.Range("AC4").Value = "Flera fruktkoppar som växer i en häxring"
$ws.Range("A5").Value = 111866276
$ws.Range("B5").Value = 78228
$ws.Range("E5").Value = 6453
$ws.Range("F5").Value = "Vedskivlav"
$ws.Range("G5").Value = "Hertelidea botryosa"
$ws.Range("H5").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q5").Value = 702661
$ws.Range("R5").Value = 7299929
$ws.Range("B6").Value = 90816
$ws.Range("A7").Value = 111865524
$ws.Range("B7").Value = 90794
$ws.Range("E7").Value = 4362
$ws.Range("F7").Value = "Blå taggsvamp"
$ws.Range("G7").Value = "Hydnellum caeruleum"
$ws.Range("H7").Value = "(Hornem.) P.Karst."
$ws.Range("Q7").Value = 702731
$ws.Range("R7").Value = 7299742
$ws.Range("B8").Value = 78228
$ws.Range("B9").Value = 90786
$ws.Range("A10").Value = 111865578
$ws.Range("B10").Value = 90988
$ws.Range("E10").Value = 2079
$ws.Range("F10").Value = "Nordtagging"
$ws.Range("G10").Value = "Odonticium romellii"
$ws.Range("H10").Value = "(S.Lundell) Parmasto"
$ws.Range("Q10").Value = 702742
$ws.Range("R10").Value = 7299746
$ws.Range("A11").Value = 111865981
$ws.Range("B11").Value = 90786
$ws.Range("E11").Value = 3100
$ws.Range("F11").Value = "Talltaggsvamp"
$ws.Range("G11").Value = "Bankera fuligineoalba"
$ws.Range("H11").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("Q11").Value = 702696
$ws.Range("R11").Value = 7299770
$ws.Range("A12").Value = 111866301
$ws.Range("B12").Value = 90794
$ws.Range("Q12").Value = 702522
$ws.Range("R12").Value = 7300048
$ws.Range("A13").Value = 111865961
$ws.Range("B13").Value = 77388
$ws.Range("E13").Value = 6446
$ws.Range("F13").Value = "Kolflarnlav"
$ws.Range("G13").Value = "Carbonicola anthracophila"
$ws.Range("H13").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q13").Value = 702714
$ws.Range("R13").Value = 7299790
$ws.Range("A14").Value = 111866021
$ws.Range("B14").Value = 78228
$ws.Range("E14").Value = 6453
$ws.Range("F14").Value = "Vedskivlav"
$ws.Range("G14").Value = "Hertelidea botryosa"
$ws.Range("H14").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q14").Value = 702738
$ws.Range("R14").Value = 7299806
$ws.Range("A15").Value = 111866170
$ws.Range("B15").Value = 90816
$ws.Range("Q15").Value = 702754
$ws.Range("R15").Value = 7299887
$ws.Range("AC15").Value = $null
$ws.Range("A16").Value = 111865488
$ws.Range("B16").Value = 90794
$ws.Range("E16").Value = 4362
$ws.Range("F16").Value = "Blå taggsvamp"
$ws.Range("G16").Value = "Hydnellum caeruleum"
$ws.Range("H16").Value = "(Hornem.) P.Karst."
$ws.Range("Q16").Value = 702716
$ws.Range("R16").Value = 7299725
$ws.Range("A17").Value = 111865919
$ws.Range("B17").Value = 95693
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 221941
$ws.Range("F17").Value = "Plattlummer"
$ws.Range("G17").Value = "Lycopodium complanatum"
$ws.Range("H17").Value = "L."
$ws.Range("Q17").Value = 702755
$ws.Range("R17").Value = 7299754
$ws.Range("A18").Value = 111866159
$ws.Range("B18").Value = 90786
$ws.Range("E18").Value = 3100
$ws.Range("F18").Value = "Talltaggsvamp"
$ws.Range("G18").Value = "Bankera fuligineoalba"
$ws.Range("H18").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("Q18").Value = 702755
$ws.Range("R18").Value = 7299865
$ws.Range("A19").Value = 111866131
$ws.Range("B19").Value = 90816
$ws.Range("E19").Value = 2059
$ws.Range("F19").Value = "Skrovlig taggsvamp"
$ws.Range("G19").Value = "Hydnellum scabrosum"
$ws.Range("H19").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q19").Value = 702757
$ws.Range("R19").Value = 7299855
$ws.Range("AC19").Value = "Flera fruktkoppar som växer i en häxring"
$ws.Range("A20").Value = 111866065
$ws.Range("B20").Value = 78228
$ws.Range("Q20").Value = 702768
$ws.Range("R20").Value = 7299828
$ws.Range("A21").Value = 111865668
$ws.Range("B21").Value = 78228
$ws.Range("E21").Value = 6453
$ws.Range("F21").Value = "Vedskivlav"
$ws.Range("G21").Value = "Hertelidea botryosa"
$ws.Range("H21").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q21").Value = 702741
$ws.Range("R21").Value = 7299744
$ws.Range("AC21").Value = $null
$ws.Range("A22").Value = 111865263
$ws.Range("B22").Value = 90792
$ws.Range("E22").Value = 4361
$ws.Range("F22").Value = "Orange taggsvamp"
$ws.Range("G22").Value = "Hydnellum aurantiacum"
$ws.Range("H22").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("R22").Value = 7299724
